$wb = $excel.ActiveWorkbook

# Add the new "Songs" worksheet after "DataTypes" (last tab) and make it active.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item("DataTypes"))
$ws.Name = "Songs"

# Header row
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Length"
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Album_Id"
$ws.Range("F1").Value = "Artist_Id"
$ws.Range("G1").Value = "Genre_Id"
$ws.Range("H1").Value = "On_Device"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Folgers Crystals"
$ws.Range("C2").Value = 2.23
$ws.Range("D2").Value = 1.29
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Orange"
$ws.Range("C3").Value = 26.03
$ws.Range("D3").Value = 9.99
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 0

# Row 4 (sparse)
$ws.Range("A4").Value = 3
$ws.Range("D4").Value = 1.29
$ws.Range("H4").Value = 0

# Turn the range into a table
$listObj = $ws.ListObjects.Add(1, $ws.Range("A1:H4"), $null, 1)
$listObj.Name = "Table3"
$listObj.TableStyle = "TableStyleMedium20"

$ws.Range("I4").Select()
$ws.Select()
